# Applies the change: fill E2:E6 with value 1, and move the active
# selection from E15 to E7 (as reflected in the sheet's sheetView).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E6").Value = 1

$ws.Range("E7").Select()
